# Applies updated excess-mortality figures (recomputed with the new
# background-population function) to the relevant country/sex rows.
#
# For each affected row we update:
#   D  - Mortality (current year deaths)
#   P  - Excess_mortality_Mean
#   R  - P_Score
#   S  - P_score_fluctuation
#   U  - Excess Mortality +/- (display string)
#   V  - P_score +/- (display string)
#   X  - Excess_mortality_per_10^5
#   Y  - Excess_mortality_per_10^5_fluc
#   Z  - Excess_mortality_per_10^5_+/- (display string)
#
# Not every row changes every one of these columns - only the ones that
# differ from the original workbook are set below.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$PM = [char]0x00B1   # "±"

function SetRowValues($Row, $Values) {
    foreach ($col in $Values.Keys) {
        $ws.Range("$col$Row").Value = $Values[$col]
    }
}

# Row 5 - Female / Croatia
SetRowValues 5 @{
    'D' = 99
    'P' = 0.6
    'R' = 0.6
    'S' = 5.4
    'U' = "0.6 (${PM}5.6)"
    'V' = "0.6% (${PM}5.4%)"
    'X' = 0.2
    'Y' = 2.2
    'Z' = "0.2(${PM}2.2)"
}

# Row 15 - Female / Italy
SetRowValues 15 @{
    'D' = 890
    'P' = -141.6
    'R' = -13.7
    'U' = "-141.6 (${PM}57.3)"
    'V' = "-13.7% (${PM}4.6%)"
}

# Row 17 - Female / Lithuania
SetRowValues 17 @{
    'D' = 128
    'P' = -19
    'R' = -12.9
    'S' = 8.2
    'U' = "-19.0 (${PM}15.3)"
    'V' = "-12.9% (${PM}8.2%)"
    'X' = -11.2
    'Z' = "-11.2(${PM}9.0)"
}

# Row 33 - Male / Bulgaria
SetRowValues 33 @{
    'D' = 794
    'P' = 62
    'R' = 8.5
    'U' = "62.0 (${PM}33.0)"
    'V' = "8.5% (${PM}4.7%)"
    'X' = 12.6
    'Z' = "12.6(${PM}6.6)"
}

# Row 34 - Male / Croatia
SetRowValues 34 @{
    'D' = 243
    'P' = -12.8
    'R' = -5
    'S' = 4.4
    'U' = "-12.8 (${PM}12.3)"
    'V' = "-5.0% (${PM}4.4%)"
    'X' = -4.6
    'Y' = 4.4
    'Z' = "-4.6(${PM}4.4)"
}

# Row 40 - Male / France
SetRowValues 40 @{
    'D' = 3416
    'P' = -34.2
    'R' = -1
    'U' = "-34.2 (${PM}65.0)"
    'V' = "-1.0% (${PM}1.8%)"
    'X' = -0.8
    'Z' = "-0.8(${PM}1.6)"
}

# Row 44 - Male / Italy
SetRowValues 44 @{
    'D' = 1769
    'P' = -175.8
    'U' = "-175.8 (${PM}119.7)"
}

# Row 46 - Male / Lithuania
SetRowValues 46 @{
    'D' = 436
    'P' = -44.2
    'R' = -9.2
    'S' = 8.5
    'U' = "-44.2 (${PM}49.8)"
    'V' = "-9.2% (${PM}8.5%)"
    'X' = -23.4
    'Z' = "-23.4(${PM}26.4)"
}

# Row 47 - Male / Luxembourg
SetRowValues 47 @{
    'D' = 26
    'P' = 1.4
    'R' = 5.7
    'S' = 7.2
    'U' = "1.4 (${PM}1.8)"
    'V' = "5.7% (${PM}7.2%)"
    'X' = 2.8
    'Z' = "2.8(${PM}3.6)"
}

# Row 53 - Male / Portugal
SetRowValues 53 @{
    'D' = 462
    'P' = -12.2
    'R' = -2.6
    'S' = 7
    'U' = "-12.2 (${PM}36.6)"
    'V' = "-2.6% (${PM}7.0%)"
    'X' = -2
    'Z' = "-2.0(${PM}6.1)"
}

# Row 62 - Total / Bulgaria
SetRowValues 62 @{
    'D' = 1103
    'P' = 58
    'R' = 5.6
    'U' = "58.0 (${PM}38.4)"
    'V' = "5.6% (${PM}3.8%)"
    'X' = 6.1
    'Z' = "6.1(${PM}4.0)"
}

# Row 63 - Total / Croatia
SetRowValues 63 @{
    'D' = 342
    'P' = -12.2
    'R' = -3.4
    'U' = "-12.2 (${PM}17.7)"
    'V' = "-3.4% (${PM}4.6%)"
    'X' = -2.3
    'Z' = "-2.3(${PM}3.3)"
}

# Row 69 - Total / France
SetRowValues 69 @{
    'D' = 4970
    'P' = -37
    'R' = -0.7
    'U' = "-37.0 (${PM}70.6)"
    'V' = "-0.7% (${PM}1.4%)"
    'X' = -0.4
    'Y' = 0.8
    'Z' = "-0.4(${PM}0.8)"
}

# Row 75 - Total / Lithuania
SetRowValues 75 @{
    'D' = 564
    'P' = -63.2
    'R' = -10.1
    'S' = 8.3
    'U' = "-63.2 (${PM}63.7)"
    'V' = "-10.1% (${PM}8.3%)"
    'X' = -17.6
    'Z' = "-17.6(${PM}17.7)"
}

# Row 76 - Total / Luxembourg
SetRowValues 76 @{
    'D' = 36
    'P' = -0.2
    'R' = -0.6
    'S' = 4.9
    'U' = "-0.2 (${PM}1.9)"
    'V' = "-0.6% (${PM}4.9%)"
    'X' = -0.2
    'Z' = "-0.2(${PM}1.9)"
}

# Row 82 - Total / Portugal
SetRowValues 82 @{
    'D' = 705
    'P' = -21.2
    'R' = -2.9
    'S' = 5
    'U' = "-21.2 (${PM}38.9)"
    'V' = "-2.9% (${PM}5.0%)"
    'X' = -1.7
    'Y' = 3.1
    'Z' = "-1.7(${PM}3.1)"
}
